$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.688.79'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '1.637.03'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  +0.00%  '
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '213.74'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("E6").Value = '  +4.37%  '
$ws.Range("E7").Value = '  +0.03%  '
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.254'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  +2.62%  '
$ws.Range("E9").Value = '  +1.58%  '
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.28'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +2.87%  '
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0845'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  +3.67%  '
$ws.Range("D12").Value = '1.865.66'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '1.635.33'
$ws.Range("E13").Value = '  +0.97%  '
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.11'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  +2.82%  '
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("D16").Value = '26.690.99'
$ws.Range("E16").Value = '  +1.59%  '
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '63.57'
$cell.Style = $origStyle
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("E18").Value = '  +2.73%  '
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '220.21'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +9.55%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  +1.15%  '
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.46'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("E24").Value = '  +1.47%  '
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '148.22'
$cell.Style = $origStyle
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  +1.50%  '
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.94'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  +6.21%  '
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.53'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("E30").Value = '  -0.51%  '
$ws.Range("E31").Value = '  +0.05%  '
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.33'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +5.14%  '
$ws.Range("E33").Value = '  +2.75%  '
$ws.Range("E34").Value = '  +2.17%  '
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("D36").Value = '1.217.25'
$ws.Range("E36").Value = '  +3.40%  '
$ws.Range("E37").Value = '  +5.95%  '
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.814'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  +1.68%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  +2.88%  '
$ws.Range("E41").Value = '  -1.07%  '
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.45'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +2.16%  '
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("D44").Value = '1.774.71'
$ws.Range("E44").Value = '  +1.13%  '
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '93.58'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("E46").Value = '  +2.56%  '
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '54.97'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("E48").Value = '  +1.09%  '
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.73'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +6.58%  '
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.410'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("E51").Value = '  +0.05%  '
